$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (shifts existing rows 3-18 down to 4-19)
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row with the new product entry
$ws.Cells.Item(3, 1).Value = "naturnes epinards des 4/6mois nestle bols 130g x 2"
$ws.Cells.Item(3, 2).Value = 2
$ws.Cells.Item(3, 3).Value = 130
$ws.Cells.Item(3, 5).Value = "g"

# Match the final selection state observed in the diff
$ws.Range("F16").Select()
